$wb = $excel.ActiveWorkbook

# Sheet ALC, row 2
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2794.4
$ws.Range("I2").Value = 3706.3572
$ws.Range("J2").Value = 666.5
$ws.Range("K2").Value = 3706.3572
$ws.Range("L2").Value = 666.5
$ws.Range("M2").Value = -3593.3572
$ws.Range("N2").Value = -892.5

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4926.2104
$ws.Range("J40").Value = 6458.5
$ws.Range("L40").Value = 6458.5
$ws.Range("N40").Value = -6808.5

# Sheet ALC, row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1094.2632
$ws.Range("I111").Value = 1094.2632
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3282.7896
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -215.7896000000001
$ws.Range("N111").ClearContents()

# Sheet ALC, row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1650.1428
$ws.Range("I131").Value = 1313.5
$ws.Range("J131").Value = 2099
$ws.Range("K131").Value = 3940.5
$ws.Range("L131").Value = 6297
$ws.Range("M131").Value = 1099.5
$ws.Range("N131").Value = -16377

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2511.7334
$ws.Range("J138").Value = 4165.3335
$ws.Range("L138").Value = 12496.0005
$ws.Range("N138").Value = -22776.0005

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19942.895
$ws.Range("I32").Value = 17520.1
$ws.Range("K32").Value = 17520.1
$ws.Range("M32").Value = -17233.1

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5753
$ws.Range("I61").Value = 5753
$ws.Range("K61").Value = 5753
$ws.Range("M61").Value = -5541

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2369.52
$ws.Range("I74").Value = 2152.2727
$ws.Range("K74").Value = 2152.2727
$ws.Range("M74").Value = -1278.2727

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2369.52
$ws.Range("I77").Value = 2152.2727
$ws.Range("K77").Value = 10761.3635
$ws.Range("M77").Value = -6393.363499999999

# Sheet ARM, row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3812.2083
$ws.Range("I102").Value = 1749.8334
$ws.Range("K102").Value = 1749.8334
$ws.Range("M102").Value = -127.8334

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5753
$ws.Range("I136").Value = 5753
$ws.Range("K136").Value = 17259
$ws.Range("M136").Value = -14709

# Sheet ARM, row 141
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 82499.5
$ws.Range("J141").Value = 82499.5
$ws.Range("L141").Value = 82499.5
$ws.Range("N141").Value = -92859.5

# Sheet BSM, row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 349
$ws.Range("I22").Value = 345.9
$ws.Range("K22").Value = 345.9
$ws.Range("M22").Value = -172.9

# Sheet BSM, row 76
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 35362.8
$ws.Range("J76").Value = 35362.8
$ws.Range("L76").Value = 35362.8
$ws.Range("N76").Value = -35992.8

# Sheet BSM, row 79
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 35362.8
$ws.Range("J79").Value = 35362.8
$ws.Range("L79").Value = 35362.8
$ws.Range("N79").Value = -37546.8

# Sheet BSM, row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2750.818
$ws.Range("I105").Value = 2161.8
$ws.Range("K105").Value = 2161.8
$ws.Range("M105").Value = -414.8000000000002

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7374.5
$ws.Range("I134").Value = 4999.5
$ws.Range("J134").Value = 9749.5
$ws.Range("K134").Value = 14998.5
$ws.Range("L134").Value = 29248.5
$ws.Range("M134").Value = -12463.5
$ws.Range("N134").Value = -34318.5

# Sheet CRP, row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 56424.832
$ws.Range("J141").Value = 56424.832
$ws.Range("L141").Value = 56424.832
$ws.Range("N141").Value = -66784.83199999999

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1915.6046
$ws.Range("I113").Value = 1361
$ws.Range("K113").Value = 4083
$ws.Range("M113").Value = -1913

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3067.9285
$ws.Range("I132").Value = 3640.4
$ws.Range("K132").Value = 32763.6
$ws.Range("M132").Value = -30233.6

# Sheet GSM, row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# Sheet GSM, row 92
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 5833.5557
$ws.Range("J92").Value = 6125.25
$ws.Range("L92").Value = 6125.25
$ws.Range("N92").Value = -9869.25

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1355.375
$ws.Range("I97").Value = 1224.6666
$ws.Range("K97").Value = 1224.6666
$ws.Range("M97").Value = -728.6666

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5290
$ws.Range("I102").Value = 5166.6665
$ws.Range("K102").Value = 5166.6665
$ws.Range("M102").Value = -3544.6665

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 10000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -14340

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 57903.5
$ws.Range("I132").Value = 66621.88
$ws.Range("K132").Value = 199865.64
$ws.Range("M132").Value = -197335.64

# Sheet LTW, row 35
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 84955.5
$ws.Range("I35").Value = 1574.5
$ws.Range("K35").Value = 1574.5
$ws.Range("M35").Value = -1238.5

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8528.166999999999
$ws.Range("J40").Value = 11393.4
$ws.Range("L40").Value = 11393.4
$ws.Range("N40").Value = -11665.4

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5069.533
$ws.Range("I61").Value = 4004.7778
$ws.Range("J61").Value = 6666.6665
$ws.Range("K61").Value = 4004.7778
$ws.Range("L61").Value = 6666.6665
$ws.Range("M61").Value = -3802.7778
$ws.Range("N61").Value = -7070.6665

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5419.7
$ws.Range("I82").Value = 4292
$ws.Range("K82").Value = 4292
$ws.Range("M82").Value = -3931

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 5419.7
$ws.Range("I85").Value = 4292
$ws.Range("K85").Value = 4292
$ws.Range("M85").Value = -3044

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5069.533
$ws.Range("I113").Value = 4004.7778
$ws.Range("J113").Value = 6666.6665
$ws.Range("K113").Value = 4004.7778
$ws.Range("L113").Value = 6666.6665
$ws.Range("M113").Value = -1834.7778
$ws.Range("N113").Value = -11006.6665

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7948.8335
$ws.Range("I136").Value = 7538.6
$ws.Range("K136").Value = 22615.8
$ws.Range("M136").Value = -20065.8

# Sheet WVR, row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 145042.86
$ws.Range("I2").Value = 145042.86
$ws.Range("K2").Value = 145042.86
$ws.Range("M2").Value = -144930.86

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11571.286
$ws.Range("I62").Value = 10499.5
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 10499.5
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -9875.5
$ws.Range("N62").Value = -13248

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 11571.286
$ws.Range("I65").Value = 10499.5
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 52497.5
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -49377.5
$ws.Range("N65").Value = -66240

# Sheet WVR, row 68
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 31333
$ws.Range("J68").Value = 31999.5
$ws.Range("L68").Value = 31999.5
$ws.Range("N68").Value = -33621.5

# Sheet WVR, row 71
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 31333
$ws.Range("J71").Value = 31999.5
$ws.Range("L71").Value = 95998.5
$ws.Range("N71").Value = -104110.5

# Sheet WVR, row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4300
$ws.Range("I81").Value = 4300
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8600
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -7539
$ws.Range("N81").ClearContents()

# Sheet WVR, row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4300
$ws.Range("I84").Value = 4300
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 43000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -37696
$ws.Range("N84").ClearContents()

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6075
$ws.Range("I132").Value = 5657.914
$ws.Range("K132").Value = 16973.742
$ws.Range("M132").Value = -14443.742

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6590.533
$ws.Range("I136").Value = 5770.7393
$ws.Range("J136").Value = 9284.143
$ws.Range("K136").Value = 17312.2179
$ws.Range("L136").Value = 27852.429
$ws.Range("M136").Value = -14762.2179
$ws.Range("N136").Value = -32952.429
